$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume values scraped on 2023-09-09
$ws.Range("D2").Value = '25.996.55'
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").Value = '1.640.83'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("E4").Value = '  +0.54%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.95'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.508'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("E8").Value = '  -0.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0638'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.59'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0795'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("D13").Value = '1.864.57'
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("D14").Value = '1.649.34'
$ws.Range("E14").Value = '  +0.24%  '
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").Value = '0.0₃0764'
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.95'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").Value = '25.887.73'
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '192.90'
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.36'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.63%  '
$ws.Range("E22").Value = '  -1.24%  '
$ws.Range("E23").Value = '  -0.39%  '
$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.64'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.41%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.80'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("E26").Value = '  +4.78%  '
$ws.Range("E27").Value = '  +0.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.92'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.54'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.41%  '
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0500'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.29'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.25'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  -3.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.47'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.902'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("D37").Value = '1.133.90'
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.543'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.67%  '
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.54'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.795'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("D44").Value = '1.773.99'
$ws.Range("E44").Value = '  -0.37%  '
$ws.Range("E45").Value = '  +2.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.62'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.42%  '
$ws.Range("E47").Value = '  +2.77%  '
$ws.Range("E48").Value = '  -0.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.73'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.49%  '
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0960'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.73%  '
